$d = $word.ActiveDocument

$replacements = @(
    @("711÷6=", "422÷7="),
    @("272÷7=", "100÷6="),
    @("528÷5=", "112÷2="),
    @("129÷6=", "311÷5="),
    @("180÷9=", "719÷7="),
    @("317÷6=", "838÷8="),
    @("532÷3=", "373÷4="),
    @("308÷4=", "853÷4="),
    @("729÷2=", "606÷2="),
    @("895÷8=", "266÷3="),
    @("237÷6=", "775÷4="),
    @("553÷2=", "128÷6="),
    @("313÷4=", "384÷2="),
    @("419÷6=", "753÷8="),
    @("479÷5=", "764÷5="),
    @("371÷5=", "979÷2="),
    @("247÷6=", "320÷7="),
    @("690÷3=", "754÷3="),
    @("962÷8=", "272÷4="),
    @("178÷2=", "547÷7="),
    @("330÷3=", "796÷5="),
    @("510÷4=", "458÷3="),
    @("410÷3=", "520÷9="),
    @("420÷2=", "875÷3="),
    @("720÷9=", "583÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
